$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Sexo" column (G) relabels its two categories:
#   Femenino  -> Mujer
#   Masculino -> Hombre
# (processed in this order so new shared-string entries land in the same
# slots the source workbook ends up with). Data rows run 2..81.
for ($r = 2; $r -le 81; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq "Femenino") {
        $cell.Value = "Mujer"
    }
}

for ($r = 2; $r -le 81; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq "Masculino") {
        $cell.Value = "Hombre"
    }
}

# Update the saved selection state on the sheet view.
$null = $ws.Range("T86").Select()
